# feat: v0.0.15 + v2.2.9 Range
# Equity curve was recomputed with a starting equity of 1000 (was 1600),
# scaling the Equity column (A) by 1000/1600 = 0.625 throughout, and
# refreshing the dependent DrawdownPct column (B) accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = 0
$ws.Range("A3").Value = 1000
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = 1000
$ws.Range("B4").Value = 0
$ws.Range("A5").Value = 1000
$ws.Range("B5").Value = 0
$ws.Range("A6").Value = 1000
$ws.Range("B6").Value = 0
$ws.Range("A7").Value = 1000
$ws.Range("B7").Value = 0
$ws.Range("A8").Value = 1000
$ws.Range("B8").Value = 0
$ws.Range("A9").Value = 1000
$ws.Range("B9").Value = 0
$ws.Range("A10").Value = 1000
$ws.Range("B10").Value = 0
$ws.Range("A11").Value = 1000
$ws.Range("B11").Value = 0
$ws.Range("A12").Value = 1000
$ws.Range("B12").Value = 0
$ws.Range("A13").Value = 1000
$ws.Range("B13").Value = 0
$ws.Range("A14").Value = 1000
$ws.Range("B14").Value = 0
$ws.Range("A15").Value = 1000
$ws.Range("B15").Value = 0
$ws.Range("A16").Value = 1000
$ws.Range("B16").Value = 0
$ws.Range("A17").Value = 1000
$ws.Range("B17").Value = 0
$ws.Range("A18").Value = 1000
$ws.Range("B18").Value = 0
$ws.Range("A19").Value = 1000
$ws.Range("B19").Value = 0
$ws.Range("A20").Value = 1000
$ws.Range("B20").Value = 0
$ws.Range("A21").Value = 1000
$ws.Range("B21").Value = 0
$ws.Range("A22").Value = 1000
$ws.Range("B22").Value = 0
$ws.Range("A23").Value = 1000
$ws.Range("B23").Value = 0
$ws.Range("A24").Value = 1000
$ws.Range("B24").Value = 0
$ws.Range("A25").Value = 1039.8018041198
$ws.Range("B25").Value = 0
$ws.Range("A26").Value = 1038.5574351198
$ws.Range("B26").Value = 0.001196736719507263
$ws.Range("A27").Value = 1017.6568567198
$ws.Range("B27").Value = 0.0212972773390655
$ws.Range("A28").Value = 1023.3152707198
$ws.Range("B28").Value = 0.01585545758304974
$ws.Range("A29").Value = 1010.3376771198
$ws.Range("B29").Value = 0.02833629147714534
$ws.Range("A30").Value = 1041.1353579198
$ws.Range("B30").Value = 0
$ws.Range("A31").Value = 1266.7289121198
$ws.Range("B31").Value = 0
$ws.Range("A32").Value = 1377.3891709198
$ws.Range("B32").Value = 0
$ws.Range("A33").Value = 1439.7233201198
$ws.Range("B33").Value = 0
$ws.Range("A34").Value = 1390.9494787198
$ws.Range("B34").Value = 0.03387723232540374
$ws.Range("A35").Value = 1321.6884425198
$ws.Range("B35").Value = 0.08198441738804241
$ws.Range("A36").Value = 1471.6285797198
$ws.Range("B36").Value = 0
$ws.Range("A37").Value = 1428.4275831198
$ws.Range("B37").Value = 0.02935591031279483
$ws.Range("A38").Value = 1360.6585845198
$ws.Range("B38").Value = 0.0754062517738876
$ws.Range("A39").Value = 1360.6585845198
$ws.Range("B39").Value = 0.0754062517738876
